# Generate Report for Handoff
# Update the localization-status workbook to reflect that the zh-cn file
# has moved from "In Translation" to "Ready for handoff", with refreshed
# handoff timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("C2").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("D2").Value = "2016-03-22 18:36:22" # Latest Handoff Date

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"    # Status
$wsZhCn.Range("E2").Value = "2016-03-22 18:36:16"  # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"    # Status
$wsDeDe.Range("E2").Value = "2016-03-22 18:36:22"  # Latest Handoff Datetime
